$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2020 Roster")

# Replace outgoing NextGen representative (YouYou Hu) with the new one (Ottavia Prunas)
$ws.Range("A22").Value = "Ottavia"
$ws.Range("B22").Value = "Prunas"
$ws.Range("D22").Value = "ottavia.prunas@unibas.ch"
$ws.Range("C22").Value = "University Hospital of Basel"
